$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.323.91'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '2.607.83'
$ws.Range("E3").Value = '  +0.68%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '543.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.07'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.47'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.67%  '
$ws.Range("E10").Value = '  +1.84%  '
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("D13").Value = '3.066.30'
$ws.Range("E13").Value = '  +0.48%  '
$ws.Range("D14").Value = '59.257.35'
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.56'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").Value = '2.617.82'
$ws.Range("E16").Value = '  +0.37%  '
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '343.92'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.37'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("E20").Value = '  -0.80%  '
$ws.Range("E21").Value = '  -1.72%  '
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("E24").Value = '  +0.64%  '
$ws.Range("E25").Value = '  -0.99%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.23'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0740'
$ws.Range("E28").Value = '  +1.65%  '
$ws.Range("B29").Value = 'USDe'
$ws.Range("C29").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.71'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.76'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.77'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.98'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '37.14'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.21%  '
$ws.Range("E36").Value = '  -1.60%  '
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.833'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.814'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.56'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.12%  '
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '277.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.598'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.12%  '
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0525'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.81%  '
$ws.Range("E47").Value = '  +1.29%  '
$ws.Range("D48").Value = '1.944.32'
$ws.Range("E48").Value = '  -2.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.41'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.51'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("E51").Value = '  -1.97%  '
